$d = $word.ActiveDocument

# Locate the target paragraphs by scanning all paragraphs for exact text
# matches, so the script is resilient to any off-by-one indexing quirks.
$ratingParaIdx = 0
$feedbackParaIdx = 0
$filmPhotographerParaIdx = 0

$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    $t = $p.Range.Text.TrimEnd("`r", "`a")
    if ($t -eq "Rating" -and $p.Style.NameLocal -eq "List Paragraph") {
        $ratingParaIdx = $idx
    }
    if ($t -eq "Feedback" -and $p.Style.NameLocal -eq "List Paragraph") {
        $feedbackParaIdx = $idx
    }
    if ($t -eq "22)Film Photographer") {
        $filmPhotographerParaIdx = $idx
    }
}

if ($ratingParaIdx -eq 0 -or $feedbackParaIdx -eq 0 -or $filmPhotographerParaIdx -eq 0) {
    throw "Could not locate one or more anchor paragraphs (Rating=$ratingParaIdx Feedback=$feedbackParaIdx Film=$filmPhotographerParaIdx)"
}

# --- Change 1: "Rating" -> "Photographer " + "Rating" (two runs, same formatting) ---
$pRating = $d.Paragraphs($ratingParaIdx)
$startRating = $pRating.Range.Start
$insRating = $d.Range($startRating, $startRating)
$insRating.InsertAfter("Photographer ")
# Toggle bold on just the inserted text to force Word to keep it as a
# separate run from the following "Rating" run, even though the final
# formatting (no bold) ends up identical on both runs.
$newRangeRating = $d.Range($startRating, $startRating + 13)
$newRangeRating.Bold = 1
$newRangeRating.Bold = 0

# --- Change 2: "Feedback" -> "Photographer " + "Feedback" (two runs, same formatting) ---
$pFeedback = $d.Paragraphs($feedbackParaIdx)
$startFeedback = $pFeedback.Range.Start
$insFeedback = $d.Range($startFeedback, $startFeedback)
$insFeedback.InsertAfter("Photographer ")
$newRangeFeedback = $d.Range($startFeedback, $startFeedback + 13)
$newRangeFeedback.Bold = 1
$newRangeFeedback.Bold = 0

# --- Change 3: add four new numbered paragraphs after "22)Film Photographer" ---
$pFilm = $d.Paragraphs($filmPhotographerParaIdx)
$pFilm.Range.InsertParagraphAfter()
$p23 = $d.Paragraphs($filmPhotographerParaIdx + 1)
$p23.Range.InsertAfter("23)Learner Feedback")

$p23.Range.InsertParagraphAfter()
$p24 = $d.Paragraphs($filmPhotographerParaIdx + 2)
$p24.Range.InsertAfter("24)Website Feedback")

$p24.Range.InsertParagraphAfter()
$p25 = $d.Paragraphs($filmPhotographerParaIdx + 3)
$p25.Range.InsertAfter("25)Learner Rating")

$p25.Range.InsertParagraphAfter()
$p26 = $d.Paragraphs($filmPhotographerParaIdx + 4)
$p26.Range.InsertAfter("26)Website Rating")

Write-Output "Rating para idx: $ratingParaIdx -> [$($d.Paragraphs($ratingParaIdx).Range.Text)]"
Write-Output "Feedback para idx: $feedbackParaIdx -> [$($d.Paragraphs($feedbackParaIdx).Range.Text)]"
Write-Output "Film Photographer para idx: $filmPhotographerParaIdx"
Write-Output "23 -> [$($p23.Range.Text)]"
Write-Output "24 -> [$($p24.Range.Text)]"
Write-Output "25 -> [$($p25.Range.Text)]"
Write-Output "26 -> [$($p26.Range.Text)]"
